$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column L that need to flip from the "X" (fail) mark to the
# "check" (pass) mark -- small bug fixed, so the "timeline" assertion now
# passes for these rows too (row 14 already showed a passing check mark).
$rows = @(13, 15, 16, 17, 18, 19, 20, 21, 22, 23)

foreach ($r in $rows) {
    $cell = $ws.Range("L$r")
    $cell.Value = [char]0x2705
    $cell.Font.Color = 32768
}
